# Adding queries in data file for Avneet's user story of BS create enrollment
# New rows 149-154 on the "SQL" worksheet (TIN enrollment / billing-service queries)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SQL")

# The queries (column B) and their descriptions (column C) were typed in first,
# row by row, followed by the row-number labels (column A) afterwards.

# ---- Row 149 : BS TIN is enrolled and Active ----
$ws.Cells.Item(149, 2).Value = "select * from OLE.BILLING_SERVICE b where b.ENRL_STS_CD='A' FETCH FIRST 1 ROW ONLY"
$ws.Cells.Item(149, 3).Value = "BS TIN is enrolled and Active"

# ---- Row 150 : BS TIN Details ----
$ws.Cells.Item(150, 2).Value = 'select * from OLE.BILLING_SERVICE b where b.IDENTIFIER_NBR=''{$tin}'' FETCH FIRST 1 ROW ONLY'
$ws.Cells.Item(150, 2).WrapText = $true
$ws.Cells.Item(150, 3).Value = "BS TIN Details"

# ---- Row 151 : BS TIN is in pending enrollment ----
$ws.Cells.Item(151, 2).Value = "select * from OLE.BILLING_SERVICE b where b.ENRL_STS_CD='PE' FETCH FIRST 1 ROW ONLY"
$ws.Cells.Item(151, 3).Value = "BS TIN is in pending enrollment"

# ---- Row 152 : BS Content Managed Validation (already enrolled) ----
$ws.Cells.Item(152, 2).Value = "select * from ole.content where TEXT_VAL like '%Your TIN/SSN is already enrolled%' order by CREAT_DTTM desc FETCH FIRST 1 ROW ONLY"
$ws.Cells.Item(152, 2).WrapText = $true
$ws.Cells.Item(152, 3).Value = "BS Content Managed Validation"

# ---- Row 153 : BS Content Managed Validation (enrolled and active) ----
$ws.Cells.Item(153, 2).Value = "select * from ole.content where TEXT_VAL like '%The TIN/SSN you entered is enrolled and active for Electronic Payments and Statements%' order by CREAT_DTTM desc FETCH FIRST 1 ROW ONLY"
$ws.Cells.Item(153, 2).WrapText = $true
$ws.Cells.Item(153, 3).Value = "BS Content Managed Validation"

# ---- Row 154 : BS Content Managed Validation (pending enrollment status) ----
$ws.Cells.Item(154, 2).Value = "select * from ole.content where TEXT_VAL like '%Your TIN is currently in a pending enrollment status%' order by CREAT_DTTM desc FETCH FIRST 1 ROW ONLY"
$ws.Cells.Item(154, 2).WrapText = $true
$ws.Cells.Item(154, 3).Value = "BS Content Managed Validation"

# Row-number labels entered last, row by row
$ws.Cells.Item(149, 1).Value = "'148"
$ws.Cells.Item(150, 1).Value = "'149"
$ws.Cells.Item(151, 1).Value = "'150"
$ws.Cells.Item(152, 1).Value = "'151"
$ws.Cells.Item(153, 1).Value = "'152"
$ws.Cells.Item(154, 1).Value = "'153"

# Leave the cursor/selection where the author left off after typing the new rows
$ws.Activate()
$ws.Cells.Item(151, 2).Select()
